# Append a new log row (row 4) to the Nalco run log sheet, matching the
# style/format of the existing data rows (e.g. row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

# Copy formatting only from row 3 onto the new row 4, so cell styles
# (alignment, etc.) match the rest of the data rows.
$ws.Range("A3:H3").Copy() | Out-Null
$ws.Range("A4:H4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Now populate the new row's values.
$ws.Cells.Item($row, 1).Value = "2025-08-12 09:41:54 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-12 15:11:54 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""
